$d = $word.ActiveDocument

# Find the first paragraph's range and update it.
$para = $d.Paragraphs(1)
$range = $para.Range

# Trim the trailing paragraph mark from the range so we only touch the text.
$range.End = $range.End - 1

# Replace text content: add two trailing spaces to the original text.
$range.Text = "This is a Microsoft word document.  "

# Now append the red-colored addition as new runs at the end of this paragraph (before the pilcrow).
$insertRange = $d.Paragraphs(1).Range
$insertRange.End = $insertRange.End - 1
$insertRange.Collapse(0)

$insertRange.InsertAfter("(This is a change – Ve")
$insertRange.Font.Color = 255
$insertRange.Collapse(0)

$insertRange.InsertAfter("rsion for main branch")
$insertRange.Font.Color = 255
$insertRange.Collapse(0)

$insertRange.InsertAfter(")")
$insertRange.Font.Color = 255
